$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Banner_Text")

# Re-worded banner captions (Maṭham -> Mutt, plus new "Acharyas" wording)
$ws3.Range("B4").Value = "Sri Brahmatantra Swatantra Parakala Swāmy Mutt Acharyas"

# New helper column holding the reusable " Maṭham " fragment
$ws3.Range("D2").Value = " Maṭham "
$ws3.Range("D4").Value = " Maṭham "
$ws3.Range("D6").Value = " Maṭham "

$ws3.Range("B6").Value = "Sri Parakāla Swāmy Mutt – The Eternal Lineage of the Sri Vedānta Deśika"
$ws3.Range("B2").Value = "Sri Parakāla Swāmy Mutt Guru Parampara"

# Allow a blank banner entry (M4 message body is now empty)
$ws3.Range("B5").ClearContents()

# Widen caption column to fit the longer text
$ws3.Columns.Item(2).ColumnWidth = 64

# Banner_Text becomes the active/visible sheet, with B3 selected
$ws3.Activate()
$ws3.Range("B3").Select()
